$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# New "disponible" (stock) values for F2:F27
$values = @(3,13,11,11,8,10,6,16,15,8,12,8,7,0,11,8,2,5,13,3,7,13,6,5,6,2)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}

# Re-format F2 with the new number format + top alignment, then propagate
# that exact format to the rest of the column via a format-only paste so
# that only a single new style entry is created (matches the authored XML).
$f2 = $ws.Cells.Item(2, 6)
$f2.NumberFormat = "#,##0"
$f2.VerticalAlignment = -4160

$f2.Copy()
$ws.Range("F3:F27").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Select cell A2, matching the saved selection in the workbook
$ws.Range("A2").Select()
